$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.076972961425781
$ws.Range("B1").Value = 3.779794216156006
$ws.Range("C1").Value = 3.675584554672241
$ws.Range("D1").Value = 3.225313186645508
$ws.Range("E1").Value = 1.249887466430664
